$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns for every row we touch so that
# numeric-looking strings (e.g. "567.87", "2.50", "0.0000231") are stored as
# text, matching the inlineStr cell type used throughout this sheet.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.902.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.90%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.895.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.81%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.87"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.44"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.77%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.501"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.82%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.894.38"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.93"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.10%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.430"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.10"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.30%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.10%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.375.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.79%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.807.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.893.80"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.80%  "

# Row 19
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.52"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.69"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.95%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.96"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.654"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.87"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.89"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.03"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.37%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.06"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -9.47%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.03"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000111"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +10.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.02"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.28%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.84%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.03"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.89%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.107"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.95%  "

# Row 34
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.60"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.03%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.954"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.56%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.39"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.18%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.90"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.40%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.84"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.91"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.115"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.61%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.13"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.50%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.11"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.51%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.269"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.81%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.694.67"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.79%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.08%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "131.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "346.21"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.40%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.02%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.17%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.65"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.12%  "
